$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 and IF headers, matching the formatting
#     used by the other header cells (bold, thin border, centered) ---
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# --- Data rows (2-13): new numeric values for columns I (I0) and J (IF) ---
$iValues = @(7, 4, 5, 7, 5, 5, 4, 1, 1, 7, 9, 8)
$jValues = @(7, 6, 7, 8, 6, 8, 5, 5, 2, 7, 9, 8)

for ($r = 2; $r -le 13; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
